$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new row 38 with the new test mail log entry.
$row = 38

$ws.Cells.Item($row, 1).Value = "Is er al nieuws?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #10: Is er al nieuws?"
$ws.Cells.Item($row, 4).Value = "Overig"
$ws.Cells.Item($row, 5).Value = "Beste afzender,`r`nBedankt voor uw e-mail. Kunt u mij meer informatie geven over waar u precies naar op zoek bent? Op basis van uw vraag kan ik nu niet direct antwoorden met het juiste nieuws. Alvast bedankt voor uw aanvullende informatie.`r`nMet vriendelijke groet,`r`n[Naam] `r`nE-mailassistent"
$ws.Cells.Item($row, 6).Value = "2025-08-03 18:33:40"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

# The multi-line "Antwoord" text auto-expands the row height; re-autofit so
# the row keeps using the sheet's default (non-custom) height, matching a
# freshly appended row rather than one that was manually resized.
$ws.Rows.Item($row).AutoFit()

# Expand the conditional formatting ranges so they cover the new row too.
$ws.Range("D2:D37").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D38"))
$ws.Range("G2:G37").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G38"))
$ws.Range("H2:H37").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H38"))
$ws.Range("I2:I37").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I38"))
$ws.Range("J2:J37").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J38"))

# Update the Dashboard summary count for the "Overig" category (10 -> 11).
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 11
